$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the wording of the B3 text value ("I do agree" -> "I agree").
$ws.Range("B3").Value = "I agree with the news story. In the homes where I was, I have experienced all the issues mentioned in the news story."

# Add a new "time" column (C) with a header and values for the first three
# data rows (row 5 is left blank, matching the source data).
$ws.Range("C1").Value = "time"
$ws.Range("C1").Font.Bold = $true

$ws.Range("C2").Value = 1.5
$ws.Range("C3").Value = 1.5
$ws.Range("C4").Value = 1.5

# Update the view/selection state to match the saved workbook.
$ws.Range("C11").Select()
